$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B127").Value = 57552
$ws.Range("E127").Value = 136.86
$ws.Range("F127").Value = -5
$ws.Range("G127").Value = -603.45
$ws.Range("B128").Value = 64329
$ws.Range("E128").Value = 128.32
$ws.Range("F128").Value = 1
$ws.Range("G128").Value = 120.69
$ws.Range("F136").Value = 18
$ws.Range("G136").Value = 1350.36
$ws.Range("B138").Value = 2407.75
$ws.Range("F149").Value = 223
$ws.Range("G149").Value = 14450.4
$ws.Range("F152").Value = 63
$ws.Range("G152").Value = 5562.27
$ws.Range("B156").Value = 30592.79
$ws.Range("F203").Value = 53
$ws.Range("G203").Value = 1068.48
$ws.Range("B216").Value = 37532.36
$ws.Range("B219").Value = 61610
$ws.Range("E219").Value = 122.71
$ws.Range("F219").Value = -58
$ws.Range("G219").Value = -5957.18
$ws.Range("B220").Value = 63565
$ws.Range("E220").Value = 109.19
$ws.Range("F220").Value = 60
$ws.Range("G220").Value = 6162.6
$ws.Range("F225").Value = 74
$ws.Range("G225").Value = 8453.02
$ws.Range("F255").Value = 529
$ws.Range("G255").Value = 90633.57000000001
$ws.Range("B260").Value = 176521.22
$ws.Range("F343").Value = 34
$ws.Range("G343").Value = 2446.98
$ws.Range("F345").Value = 40
$ws.Range("G345").Value = 2456.4
$ws.Range("B346").Value = 24747.51
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("B463").Value = 64833
$ws.Range("E463").Value = 34.9
$ws.Range("F463").Value = 95
$ws.Range("G463").Value = 3118.85
$ws.Range("B464").Value = 60025
$ws.Range("E464").Value = 37.22
$ws.Range("F464").Value = -98
$ws.Range("G464").Value = -3217.34
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 107
$ws.Range("G474").Value = 3512.81
$ws.Range("F539").Value = 24
$ws.Range("G539").Value = 6217.68
$ws.Range("B547").Value = 16430.36
$ws.Range("F599").Value = 1461
$ws.Range("G599").Value = 238303.71
$ws.Range("F601").Value = 381
$ws.Range("G601").Value = 107773.47
$ws.Range("F602").Value = 319
$ws.Range("G602").Value = 46143.35
$ws.Range("B606").Value = 393068.58
$ws.Range("F613").Value = 132
$ws.Range("G613").Value = 21009.12
$ws.Range("B618").Value = 42559.32
$ws.Range("B619").Value = 1640526.93
$ws.Range("B620").Value = 1640526.93
